# Resolved issues with admin erasing bay info
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fix corrupted row 2 (admin accidentally erased/overwrote bay info) ---
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = "teet"
$ws.Range("D2").Value = "t4etw"

# --- Re-add the missing bay rows (44-46) ---
$ws.Range("A44").Value = 1
$ws.Range("B44").Value = 1
$ws.Range("C44").Value = "Test1"
$ws.Range("D44").Value = "Test1"
$ws.Range("E44").Value = 1
$ws.Range("F44").Value = "https://www.youtube.com/embed/Ov1v-PxiFMU?autoplay=1&mute=0|https://www.youtube.com/embed/gsIQjyeBC_c?autoplay=1&mute=0"

$ws.Range("A45").Value = 1
$ws.Range("B45").Value = 2
$ws.Range("C45").Value = "Test2dwadadw"
$ws.Range("D45").Value = "Test2"
$ws.Range("E45").Value = 1
$ws.Range("F45").Value = "https://www.youtube.com/embed/gsIQjyeBC_c?autoplay=1&mute=0|https://www.youtube.com/embed/dQw4w9WgXcQ?autoplay=1&mute=0"

$ws.Range("A46").Value = 2
$ws.Range("B46").Value = 1
$ws.Range("C46").Value = "test2"
$ws.Range("D46").Value = "test2"
$ws.Range("E46").Value = 2
$ws.Range("F46").Value = "https://example.com/metrics3|https://example.com/metrics4"
